# Weekly update: insert a new Fruta/Limón price record at row 227,
# shifting the existing rows 227-244 down to 228-245.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 227 (pushes old 227..244 -> 228..245).
$ws.Rows("227").Insert()

# Populate the new row 227 with the new weekly record.
$ws.Cells.Item(227, 1).Value  = 1
$ws.Cells.Item(227, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(227, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(227, 4).Value  = 44714
$ws.Cells.Item(227, 5).Value  = 15
$ws.Cells.Item(227, 6).Value  = "Fruta"
$ws.Cells.Item(227, 7).Value  = 100102
$ws.Cells.Item(227, 8).Value  = "Cítricos"
$ws.Cells.Item(227, 9).Value  = 100102003
$ws.Cells.Item(227, 10).Value = "Limón"
$ws.Cells.Item(227, 11).Value = "Sin especificar"
$ws.Cells.Item(227, 12).Value = "2a amarillo"
$ws.Cells.Item(227, 13).Value = 300
$ws.Cells.Item(227, 14).Value = 15000
$ws.Cells.Item(227, 15).Value = 16000
$ws.Cells.Item(227, 16).Value = 15500
$ws.Cells.Item(227, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(227, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(227, 19).Value = 775
$ws.Cells.Item(227, 20).Value = 20
